$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hungary NB I")

# --- Row 132 ---
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 1
$ws.Range("J132").Value = "A"
$ws.Range("R132").Value = 2
$ws.Range("S132").Value = 1.85
$ws.Range("W132").Value = -1
$ws.Range("X132").Value = -1
$ws.Range("Y132").Value = 5
$ws.Range("Z132").Value = -1
$ws.Range("AA132").Value = 0.8500000000000001
$ws.Range("AB132").Value = -1
$ws.Range("AC132").Value = 0.875

# --- Row 133 ---
$ws.Range("H133").Value = 1
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = "H"
$ws.Range("R133").Value = 1.8
$ws.Range("S133").Value = 2.05
$ws.Range("U133").Value = 1.9
$ws.Range("V133").Value = 1.95
$ws.Range("W133").Value = 1
$ws.Range("X133").Value = -1
$ws.Range("Y133").Value = -1
$ws.Range("Z133").Value = 0.8
$ws.Range("AA133").Value = -1
$ws.Range("AB133").Value = -1
$ws.Range("AC133").Value = 0.95

# --- Row 134 ---
$ws.Range("H134").Value = 1
$ws.Range("I134").Value = 1
$ws.Range("J134").Value = "D"
$ws.Range("U134").Value = 1.85
$ws.Range("V134").Value = 2
$ws.Range("W134").Value = -1
$ws.Range("X134").Value = 2.25
$ws.Range("Y134").Value = -1
$ws.Range("Z134").Value = 0.4125
$ws.Range("AA134").Value = -0.5
$ws.Range("AB134").Value = -1
$ws.Range("AC134").Value = 1

# --- Row 135 ---
$ws.Range("N135").Value = 3.3
$ws.Range("U135").Value = 1.85
$ws.Range("V135").Value = 2

# --- Row 136 ---
$ws.Range("N136").Value = 1.333
$ws.Range("O136").Value = 5.75
$ws.Range("P136").Value = 6.5
$ws.Range("R136").Value = 1.95
$ws.Range("S136").Value = 1.9
$ws.Range("U136").Value = 2.025
$ws.Range("V136").Value = 1.825

# --- Row 137 ---
$ws.Range("N137").Value = 2.375
$ws.Range("P137").Value = 2.75
$ws.Range("Q137").Value = 0
$ws.Range("R137").Value = 1.775
$ws.Range("S137").Value = 2.1
$ws.Range("U137").Value = 1.95
$ws.Range("V137").Value = 1.9
